$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This workbook had batting/bowling scorecard data duplicated from another
# match; this script rewrites the affected cells with the correct figures
# for this Sunrisers Sri Lanka vs Kolkata England match.

# --- Batting scorecards (rows 2-12) ---
# Row 2
$ws.Range("A2").Value = "Jason Roy"
$ws.Range("B2").Value = 2
$ws.Range("C2").Value = 2
$ws.Range("D2").Value = "LBW"
$ws.Range("E2").Value = " Nuwan Pradeep"
$ws.Range("F2").Value = 1
$ws.Range("J2").Value = "Dinesh Chandimal"
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 2
$ws.Range("M2").Value = "LBW"
$ws.Range("N2").Value = " Mark Wood"
$ws.Range("O2").Value = 1

# Row 3
$ws.Range("A3").Value = "Jos Buttler"
$ws.Range("B3").Value = 41
$ws.Range("C3").Value = 18
$ws.Range("D3").Value = "Caught"
$ws.Range("E3").Value = " Chamika Karunarathne"
$ws.Range("F3").Value = 2
$ws.Range("J3").Value = "Pathum Nissanka"
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 2
$ws.Range("M3").Value = "Caught"
$ws.Range("N3").Value = " Chris Jordan"
$ws.Range("O3").Value = 2

# Row 4
$ws.Range("A4").Value = "Dawid Malan"
$ws.Range("B4").Value = 0
$ws.Range("C4").Value = 3
$ws.Range("D4").Value = "LBW"
$ws.Range("E4").Value = " Nuwan Pradeep"
$ws.Range("F4").Value = 3
$ws.Range("J4").Value = "Charith Asalanka"
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = "Bowled"
$ws.Range("N4").Value = " Mark Wood"
$ws.Range("O4").Value = 3

# Row 5
$ws.Range("A5").Value = "Jonny Bairstow"
$ws.Range("B5").Value = 32
$ws.Range("C5").Value = 10
$ws.Range("D5").Value = "Bowled"
$ws.Range("E5").Value = " Chamika Karunarathne"
$ws.Range("F5").Value = 4
$ws.Range("J5").Value = "Dhananjaya de Silva"
$ws.Range("K5").Value = 14
$ws.Range("L5").Value = 5
$ws.Range("M5").Value = "Bowled"
$ws.Range("N5").Value = " Adil Rashid"
$ws.Range("O5").Value = 4

# Row 6
$ws.Range("A6").Value = "Eoin Morgan(C)"
$ws.Range("B6").Value = 2
$ws.Range("C6").Value = 2
$ws.Range("D6").Value = "Caught"
$ws.Range("E6").Value = " Chamika Karunarathne"
$ws.Range("F6").Value = 5
$ws.Range("J6").Value = "Bhanuka Rajapakse"
$ws.Range("K6").Value = 7
$ws.Range("L6").Value = 4
$ws.Range("M6").Value = "Bowled"
$ws.Range("N6").Value = " Chris Jordan"
$ws.Range("O6").Value = 5

# Row 7
$ws.Range("A7").Value = "Moeen Ali"
$ws.Range("B7").Value = 34
$ws.Range("C7").Value = 11
$ws.Range("D7").Value = "LBW"
$ws.Range("E7").Value = " Nuwan Pradeep"
$ws.Range("F7").Value = 6
$ws.Range("J7").Value = "Dasun Shanka(C)"
$ws.Range("K7").Value = 37
$ws.Range("L7").Value = 17
$ws.Range("M7").Value = "Bowled"
$ws.Range("N7").Value = " Mark Wood"
$ws.Range("O7").Value = 6

# Row 8
$ws.Range("A8").Value = "Liam Livingstone"
$ws.Range("B8").Value = 21
$ws.Range("C8").Value = 8
$ws.Range("D8").Value = "NOT OUT"
$ws.Range("E8").Value = " "
$ws.Range("F8").Value = 7
$ws.Range("J8").Value = "Wanindu Hasaranga"
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = "LBW"
$ws.Range("N8").Value = " Chris Jordan"
$ws.Range("O8").Value = 7

# Row 9
$ws.Range("A9").Value = "Chris Woakes"
$ws.Range("B9").Value = 5
$ws.Range("C9").Value = 2
$ws.Range("D9").Value = "Bowled"
$ws.Range("E9").Value = " Wanindu Hasaranga"
$ws.Range("F9").Value = 8
$ws.Range("J9").Value = "Chamika Karunarathne"
$ws.Range("K9").Value = 15
$ws.Range("L9").Value = 6
$ws.Range("M9").Value = "LBW"
$ws.Range("N9").Value = " Chris Woakes"
$ws.Range("O9").Value = 8

# Row 10
$ws.Range("A10").Value = "Chris Jordan"
$ws.Range("B10").Value = 9
$ws.Range("C10").Value = 3
$ws.Range("D10").Value = "Bowled"
$ws.Range("E10").Value = " Wanindu Hasaranga"
$ws.Range("F10").Value = 9
$ws.Range("J10").Value = "Dushmantha Chameera"
$ws.Range("K10").Value = 51
$ws.Range("L10").Value = 22
$ws.Range("M10").Value = "Caught"
$ws.Range("N10").Value = " Liam Livingstone"
$ws.Range("O10").Value = 9

# Row 11
$ws.Range("A11").Value = "Adil Rashid"
$ws.Range("B11").Value = 12
$ws.Range("C11").Value = 4
$ws.Range("D11").Value = "LBW"
$ws.Range("E11").Value = " Nuwan Pradeep"
$ws.Range("F11").Value = 10
$ws.Range("J11").Value = "Maheesh Theekshana"
$ws.Range("K11").Value = 6
$ws.Range("L11").Value = 3
$ws.Range("M11").Value = "NOT OUT"
$ws.Range("N11").Value = " "
$ws.Range("O11").Value = 10

# Row 12
$ws.Range("A12").Value = "Mark Wood"
$ws.Range("B12").Value = 6
$ws.Range("C12").Value = 2
$ws.Range("D12").Value = "LBW"
$ws.Range("E12").Value = " Nuwan Pradeep"
$ws.Range("F12").Value = 11
$ws.Range("J12").Value = "Nuwan Pradeep"
$ws.Range("K12").Value = 4
$ws.Range("L12").Value = 2
$ws.Range("M12").Value = "Bowled"
$ws.Range("N12").Value = " Mark Wood"
$ws.Range("O12").Value = 11

# --- Innings totals (row 16) ---
$ws.Range("A16").Value = 164
$ws.Range("B16").Value = 10
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = "10.5"
$ws.Range("C16").ClearFormats()
$ws.Range("D16").Value = 65
$ws.Range("J16").Value = 139
$ws.Range("K16").Value = 10
$ws.Range("L16").NumberFormat = "@"
$ws.Range("L16").Value = "10.5"
$ws.Range("L16").ClearFormats()
$ws.Range("M16").Value = 65

# --- Bowling figures (rows 21-25) ---
# Row 21
$ws.Range("A21").Value = "Maheesh Theekshana"
$ws.Range("B21").NumberFormat = "@"
$ws.Range("B21").Value = "2.0"
$ws.Range("B21").ClearFormats()
$ws.Range("C21").Value = 39
$ws.Range("D21").Value = 0
$ws.Range("E21").Value = 19.5
$ws.Range("J21").Value = "Adil Rashid"
$ws.Range("K21").NumberFormat = "@"
$ws.Range("K21").Value = "2.0"
$ws.Range("K21").ClearFormats()
$ws.Range("L21").Value = 23
$ws.Range("M21").Value = 1
$ws.Range("N21").Value = 11.5

# Row 22
$ws.Range("A22").Value = "Dushmantha Chameera"
$ws.Range("B22").NumberFormat = "@"
$ws.Range("B22").Value = "2.0"
$ws.Range("B22").ClearFormats()
$ws.Range("C22").Value = 24
$ws.Range("D22").Value = 0
$ws.Range("E22").Value = 12
$ws.Range("J22").Value = "Chris Jordan"
$ws.Range("K22").NumberFormat = "@"
$ws.Range("K22").Value = "2.0"
$ws.Range("K22").ClearFormats()
$ws.Range("L22").Value = 24
$ws.Range("M22").Value = 3
$ws.Range("N22").Value = 12

# Row 23
$ws.Range("A23").Value = "Chamika Karunarathne"
$ws.Range("B23").NumberFormat = "@"
$ws.Range("B23").Value = "2.0"
$ws.Range("B23").ClearFormats()
$ws.Range("C23").Value = 29
$ws.Range("D23").Value = 3
$ws.Range("E23").Value = 14.5
$ws.Range("J23").Value = "Chris Woakes"
$ws.Range("K23").NumberFormat = "@"
$ws.Range("K23").Value = "2.0"
$ws.Range("K23").ClearFormats()
$ws.Range("L23").Value = 27
$ws.Range("M23").Value = 1
$ws.Range("N23").Value = 13.5

# Row 24
$ws.Range("A24").Value = "Wanindu Hasaranga"
$ws.Range("B24").NumberFormat = "@"
$ws.Range("B24").Value = "2.0"
$ws.Range("B24").ClearFormats()
$ws.Range("C24").Value = 39
$ws.Range("D24").Value = 2
$ws.Range("E24").Value = 19.5
$ws.Range("J24").Value = "Liam Livingstone"
$ws.Range("K24").NumberFormat = "@"
$ws.Range("K24").Value = "2.0"
$ws.Range("K24").ClearFormats()
$ws.Range("L24").Value = 33
$ws.Range("M24").Value = 1
$ws.Range("N24").Value = 16.5

# Row 25
$ws.Range("A25").Value = "Nuwan Pradeep"
$ws.Range("B25").NumberFormat = "@"
$ws.Range("B25").Value = "2.5"
$ws.Range("B25").ClearFormats()
$ws.Range("C25").Value = 33
$ws.Range("D25").Value = 5
$ws.Range("E25").Value = 13.2
$ws.Range("J25").Value = "Mark Wood"
$ws.Range("K25").NumberFormat = "@"
$ws.Range("K25").Value = "2.5"
$ws.Range("K25").ClearFormats()
$ws.Range("L25").Value = 32
$ws.Range("M25").Value = 4
$ws.Range("N25").Value = 12.8

